$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column H: a comma-joined "," & A & "," & C helper column used to
# collect remote codes from a new receive sketch. Enter the first cell,
# then fill the rest with two separate range-formula assignments so the
# resulting file keeps two shared-formula groups (H4:H67 and H68:H70),
# matching how the data was originally keyed in by the author (first
# typed at H3, filled down in two passes).
$ws.Range("H3").Formula = '=","&A3&","&C3'
$ws.Range("H4:H67").Formula = '=","&A4&","&C4'
$ws.Range("H68:H70").Formula = '=","&A68&","&C68'

# Update the active selection to reflect where the author ended up after
# keying in the new column.
$ws.Range("I35").Select() | Out-Null
